# Updated cryptos list on Wed Jan 10 15:44:23 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "45.321.78"
$ws.Range("E2").Value = "  -3.42%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.442.51"
$ws.Range("E3").Value = "  +7.80%  "

# Row 5 - BNB
$ws.Range("D5").Value = "294.09"
$ws.Range("E5").Value = "  -2.25%  "

# Row 6 - Solana
$ws.Range("D6").Value = "93.59"
$ws.Range("E6").Value = "  -6.57%  "

# Row 7 - XRP
$ws.Range("D7").Value = "0.559"
$ws.Range("E7").Value = "  -0.41%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.13%  "

# Row 9 - Cardano
$ws.Range("D9").Value = "0.502"
$ws.Range("E9").Value = "  -1.96%  "

# Row 10 - Avalanche
$ws.Range("D10").Value = "34.27"
$ws.Range("E10").Value = "  -3.83%  "

# Row 11 - Dogecoin
$ws.Range("E11").Value = "  -0.61%  "

# Row 12 - Polkadot
$ws.Range("D12").Value = "7.02"
$ws.Range("E12").Value = "  -2.18%  "

# Row 13 - TRON
$ws.Range("E13").Value = "  +1.72%  "

# Row 14 - WrappedliquidstakedEther2.0
$ws.Range("D14").Value = "2.817.74"
$ws.Range("E14").Value = "  +7.96%  "

# Row 15 - WrappedEther
$ws.Range("D15").Value = "2.429.10"
$ws.Range("E15").Value = "  +7.28%  "

# Row 16 - Chainlink
$ws.Range("D16").Value = "14.26"
$ws.Range("E16").Value = "  +5.00%  "

# Row 17 - Polygon
$ws.Range("D17").Value = "0.840"
$ws.Range("E17").Value = "  +5.66%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "45.355.22"
$ws.Range("E18").Value = "  -3.19%  "

# Row 19 - InternetComputer(DFINITY)
$ws.Range("E19").Value = "  -2.99%  "

# Row 20 - ShibaInu
$ws.Range("D20").Value = "0.0₃0937"
$ws.Range("E20").Value = "  +1.20%  "

# Row 21 - Uniswap
$ws.Range("D21").Value = "6.21"
$ws.Range("E21").Value = "  +6.00%  "

# Row 22 - Litecoin
$ws.Range("D22").Value = "66.95"
$ws.Range("E22").Value = "  +2.84%  "

# Row 23 - BitcoinCash
$ws.Range("D23").Value = "239.25"
$ws.Range("E23").Value = "  -3.89%  "

# Row 24 - PancakeSwap
$ws.Range("D24").Value = "2.77"
$ws.Range("E24").Value = "  -1.45%  "

# Row 25 - Dai
$ws.Range("E25").Value = "  +0.16%  "

# Row 26 - ImmutableX
$ws.Range("D26").Value = "1.92"
$ws.Range("E26").Value = "  +2.77%  "

# Row 27
$ws.Range("E27").Value = "  -1.00%  "

# Row 28
$ws.Range("D28").Value = "37.23"
$ws.Range("E28").Value = "  -11.96%  "

# Row 29
$ws.Range("D29").Value = "9.59"
$ws.Range("E29").Value = "  -1.19%  "

# Row 30
$ws.Range("D30").Value = "3.89"
$ws.Range("E30").Value = "  +22.07%  "

# Row 31
$ws.Range("D31").Value = "21.44"
$ws.Range("E31").Value = "  +7.92%  "

# Row 32
$ws.Range("D32").Value = "149.17"
$ws.Range("E32").Value = "  +2.64%  "

# Row 33
$ws.Range("D33").Value = "2.73"
$ws.Range("E33").Value = "  -2.10%  "

# Row 34
$ws.Range("E34").Value = "  -0.22%  "

# Row 35 / Row 36 - Hedera and ARBITRUM swap places
$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").Value = "2.02"
$ws.Range("E35").Value = "  +17.95%  "

$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").Value = "0.0763"
$ws.Range("E36").Value = "  -1.46%  "

# Row 37
$ws.Range("E37").Value = "  -2.01%  "

# Row 38
$ws.Range("E38").Value = "  -0.53%  "

# Row 39
$ws.Range("D39").Value = "14.37"
$ws.Range("E39").Value = "  -11.60%  "

# Row 40
$ws.Range("D40").Value = "3.73"
$ws.Range("E40").Value = "  -1.98%  "

# Row 41
$ws.Range("E41").Value = "  -1.40%  "

# Row 42
$ws.Range("D42").Value = "1.999.39"
$ws.Range("E42").Value = "  +12.14%  "

# Row 43
$ws.Range("D43").Value = "3.17"
$ws.Range("E43").Value = "  -1.21%  "

# Row 44
$ws.Range("D44").Value = "0.998"
$ws.Range("E44").Value = "  -0.10%  "

# Row 45
$ws.Range("D45").Value = "88.24"
$ws.Range("E45").Value = "  -3.52%  "

# Row 46
$ws.Range("D46").Value = "16.22"
$ws.Range("E46").Value = "  +26.11%  "

# Row 47
$ws.Range("E47").Value = "  -13.72%  "

# Row 48
$ws.Range("D48").Value = "8.60"
$ws.Range("E48").Value = "  +9.87%  "

# Row 49
$ws.Range("D49").Value = "101.79"
$ws.Range("E49").Value = "  +8.42%  "

# Row 50
$ws.Range("D50").Value = "2.680.04"
$ws.Range("E50").Value = "  +7.74%  "

# Row 51
$ws.Range("E51").Value = "  -3.59%  "
